# Scheduled data refresh for Unicorn_Profits workbook.
# Updates cached market-board derived figures (currentAveragePrice /
# currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ / LevePriceHQ /
# LeveProfitNQ / LeveProfitHQ) for the leves whose item prices moved since
# the last scheduled run, across the ALC/ARM/BSM/CRP/CUL/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet ALC: 28 cell(s) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 33338544
$ws.Range("I106").Value = 333333340
$ws.Range("J106").Value = 5788.8887
$ws.Range("K106").Value = 333333340
$ws.Range("L106").Value = 5788.8887
$ws.Range("M106").Value = -333332709
$ws.Range("N106").Value = -7050.8887
$ws.Range("H111").Value = 1280
$ws.Range("I111").Value = 472.66666
$ws.Range("J111").Value = 1764.4
$ws.Range("K111").Value = 1417.99998
$ws.Range("L111").Value = 5293.200000000001
$ws.Range("M111").Value = 1649.00002
$ws.Range("N111").Value = -11427.2
$ws.Range("H132").Value = 2657.3403
$ws.Range("I132").Value = 889.9756
$ws.Range("J132").Value = 14734.333
$ws.Range("K132").Value = 2669.9268
$ws.Range("L132").Value = 44202.999
$ws.Range("M132").Value = -139.9268000000002
$ws.Range("N132").Value = -49262.999
$ws.Range("H137").Value = 2283379.5
$ws.Range("I137").Value = 2570937.2
$ws.Range("J137").Value = 1909554.5
$ws.Range("K137").Value = 7712811.600000001
$ws.Range("L137").Value = 5728663.5
$ws.Range("M137").Value = -7710261.600000001
$ws.Range("N137").Value = -5733763.5

# --- Sheet ARM: 39 cell(s) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 920.53845
$ws.Range("I2").Value = 1062
$ws.Range("J2").Value = 449
$ws.Range("K2").Value = 1062
$ws.Range("L2").Value = 449
$ws.Range("M2").Value = -949
$ws.Range("N2").Value = -675
$ws.Range("H45").Value = 715.2222
$ws.Range("I45").Value = 679.625
$ws.Range("K45").Value = 679.625
$ws.Range("M45").Value = -302.625
$ws.Range("H74").Value = 10589560
$ws.Range("I74").Value = 6358679
$ws.Range("J74").Value = 33436318
$ws.Range("K74").Value = 6358679
$ws.Range("L74").Value = 33436318
$ws.Range("M74").Value = -6357805
$ws.Range("N74").Value = -33438066
$ws.Range("H77").Value = 10589560
$ws.Range("I77").Value = 6358679
$ws.Range("J77").Value = 33436318
$ws.Range("K77").Value = 31793395
$ws.Range("L77").Value = 167181590
$ws.Range("M77").Value = -31789027
$ws.Range("N77").Value = -167190326
$ws.Range("H116").Value = 920.53845
$ws.Range("I116").Value = 1062
$ws.Range("J116").Value = 449
$ws.Range("K116").Value = 1062
$ws.Range("L116").Value = 449
$ws.Range("M116").Value = 1232
$ws.Range("N116").Value = -5037
$ws.Range("H122").Value = 2779.2354
$ws.Range("I122").Value = 2431.3333
$ws.Range("J122").Value = 5388.5
$ws.Range("K122").Value = 7293.999899999999
$ws.Range("L122").Value = 16165.5
$ws.Range("M122").Value = -4843.999899999999
$ws.Range("N122").Value = -21065.5

# --- Sheet BSM: 36 cell(s) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 920.53845
$ws.Range("I3").Value = 1062
$ws.Range("J3").Value = 449
$ws.Range("K3").Value = 1062
$ws.Range("L3").Value = 449
$ws.Range("M3").Value = -948
$ws.Range("N3").Value = -677
$ws.Range("H80").Value = 343.3
$ws.Range("J80").Value = 295
$ws.Range("L80").Value = 295
$ws.Range("N80").Value = -2291
$ws.Range("H83").Value = 343.3
$ws.Range("J83").Value = 295
$ws.Range("L83").Value = 1475
$ws.Range("N83").Value = -11459
$ws.Range("H86").Value = 4032.3704
$ws.Range("I86").Value = 5174.9375
$ws.Range("J86").Value = 2370.4546
$ws.Range("K86").Value = 5174.9375
$ws.Range("L86").Value = 2370.4546
$ws.Range("M86").Value = -4051.9375
$ws.Range("N86").Value = -4616.4546
$ws.Range("H89").Value = 4032.3704
$ws.Range("I89").Value = 5174.9375
$ws.Range("J89").Value = 2370.4546
$ws.Range("K89").Value = 25874.6875
$ws.Range("L89").Value = 11852.273
$ws.Range("M89").Value = -20258.6875
$ws.Range("N89").Value = -23084.273
$ws.Range("H107").Value = 1003.3461
$ws.Range("I107").Value = 982.7917
$ws.Range("J107").Value = 1250
$ws.Range("K107").Value = 982.7917
$ws.Range("L107").Value = 1250
$ws.Range("M107").Value = 937.2083
$ws.Range("N107").Value = -5090

# --- Sheet CRP: 42 cell(s) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1718693.2
$ws.Range("I31").Value = 1055.5652
$ws.Range("J31").Value = 3035548.8
$ws.Range("K31").Value = 1055.5652
$ws.Range("L31").Value = 3035548.8
$ws.Range("M31").Value = -760.5652
$ws.Range("N31").Value = -3036138.8
$ws.Range("H34").Value = 1718693.2
$ws.Range("I34").Value = 1055.5652
$ws.Range("J34").Value = 3035548.8
$ws.Range("K34").Value = 1055.5652
$ws.Range("L34").Value = 3035548.8
$ws.Range("M34").Value = -853.5652
$ws.Range("N34").Value = -3035952.8
$ws.Range("H58").Value = 4275.244
$ws.Range("I58").Value = 5694.25
$ws.Range("J58").Value = 2923.8096
$ws.Range("K58").Value = 5694.25
$ws.Range("L58").Value = 2923.8096
$ws.Range("M58").Value = -5491.25
$ws.Range("N58").Value = -3329.8096
$ws.Range("H132").Value = 10640185
$ws.Range("I132").Value = 19232018
$ws.Range("J132").Value = 2678.762
$ws.Range("K132").Value = 57696054
$ws.Range("L132").Value = 8036.286
$ws.Range("M132").Value = -57693524
$ws.Range("N132").Value = -13096.286
$ws.Range("H134").Value = 15626306
$ws.Range("I134").Value = 20000710
$ws.Range("J134").Value = 3434.8572
$ws.Range("K134").Value = 60002130
$ws.Range("L134").Value = 10304.5716
$ws.Range("M134").Value = -59999595
$ws.Range("N134").Value = -15374.5716
$ws.Range("H136").Value = 4275.244
$ws.Range("I136").Value = 5694.25
$ws.Range("J136").Value = 2923.8096
$ws.Range("K136").Value = 17082.75
$ws.Range("L136").Value = 8771.4288
$ws.Range("M136").Value = -14532.75
$ws.Range("N136").Value = -13871.4288

# --- Sheet CUL: 46 cell(s) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 95682.91
$ws.Range("I70").Value = 251628
$ws.Range("J70").Value = 6571.4287
$ws.Range("K70").Value = 754884
$ws.Range("L70").Value = 19714.2861
$ws.Range("M70").Value = -754569
$ws.Range("N70").Value = -20344.2861
$ws.Range("H73").Value = 95682.91
$ws.Range("I73").Value = 251628
$ws.Range("J73").Value = 6571.4287
$ws.Range("K73").Value = 754884
$ws.Range("L73").Value = 19714.2861
$ws.Range("M73").Value = -753792
$ws.Range("N73").Value = -21898.2861
$ws.Range("H131").Value = 1226.2543
$ws.Range("J131").Value = 1201.38
$ws.Range("L131").Value = 3604.14
$ws.Range("N131").Value = -13684.14
$ws.Range("H132").Value = 1378.1143
$ws.Range("I132").Value = 2140.1538
$ws.Range("J132").Value = 927.8182
$ws.Range("K132").Value = 19261.3842
$ws.Range("L132").Value = 8350.363800000001
$ws.Range("M132").Value = -16731.3842
$ws.Range("N132").Value = -13410.3638
$ws.Range("H136").Value = 2937.5
$ws.Range("I136").Value = 2416.6667
$ws.Range("J136").Value = 4500
$ws.Range("K136").Value = 7250.000100000001
$ws.Range("L136").Value = 13500
$ws.Range("M136").Value = -2150.000100000001
$ws.Range("N136").Value = -23700
$ws.Range("H137").Value = 8110.3125
$ws.Range("I137").Value = 1414.2858
$ws.Range("J137").Value = 13318.333
$ws.Range("K137").Value = 4242.857400000001
$ws.Range("L137").Value = 39954.999
$ws.Range("M137").Value = 857.1425999999992
$ws.Range("N137").Value = -50154.999
$ws.Range("H141").Value = 2724.389
$ws.Range("I141").Value = 2289.2144
$ws.Range("J141").Value = 4247.5
$ws.Range("K141").Value = 6867.6432
$ws.Range("L141").Value = 12742.5
$ws.Range("M141").Value = -1687.6432
$ws.Range("N141").Value = -23102.5

# --- Sheet LTW: 18 cell(s) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1440.3684
$ws.Range("I16").Value = 1440.3684
$ws.Range("K16").Value = 1440.3684
$ws.Range("M16").Value = -1270.3684
$ws.Range("H55").Value = 177.05882
$ws.Range("I55").Value = 154.66667
$ws.Range("J55").Value = 202.25
$ws.Range("K55").Value = 154.66667
$ws.Range("L55").Value = 202.25
$ws.Range("M55").Value = 18.33332999999999
$ws.Range("N55").Value = -548.25
$ws.Range("H132").Value = 15163479
$ws.Range("I132").Value = 3573.0908
$ws.Range("J132").Value = 30323386
$ws.Range("K132").Value = 10719.2724
$ws.Range("L132").Value = 90970158
$ws.Range("M132").Value = -8189.2724
$ws.Range("N132").Value = -90975218

# --- Sheet WVR: 7 cell(s) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 15805791
$ws.Range("I136").Value = 17433928
$ws.Range("J136").Value = 8541793
$ws.Range("K136").Value = 52301784
$ws.Range("L136").Value = 25625379
$ws.Range("M136").Value = -52299234
$ws.Range("N136").Value = -25630479

Write-Output "Updated $([int]216) cell(s) across 7 sheets (ALC, ARM, BSM, CRP, CUL, LTW, WVR)."
